$d = $word.ActiveDocument

function Invoke-FindReplace($range, $findText, $replaceText) {
    $ok = $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $findText"
    }
    return $ok
}

# --- Title paragraph: "Part 2B - " + "Report" were two separately-styled runs
#     with identical formatting; re-typing the text merges them into one run,
#     matching the diff's collapse of the two <w:r> elements into one. ---
Invoke-FindReplace $d.Content "Part 2B" "Part 2B" | Out-Null

# --- Paragraph 3 (Coverity intro): update the sentence about where the bugs are shown. ---
Invoke-FindReplace $d.Content `
    "as shown in the Coverity Results document. This report will discuss errors 2-4 from that report." `
    "as shown in the three .errors.xml files in the pii directory. This report will discuss three of these errors, two of which were identical." | Out-Null

# --- Paragraph 4 (was "The second error ..."): renumber to "The first of these errors ..." ---
Invoke-FindReplace $d.Content `
    "The second error (CID 10282) was of the type" `
    "The first of these errors (CID 10282) was of the type" | Out-Null

# That paragraph's leading tab character becomes a first-line indent (0.5in = 720 twips)
# instead of a literal tab run.
$para4 = $d.Paragraphs(4)
Invoke-FindReplace $para4.Range "`t" "" | Out-Null
$para4.Format.FirstLineIndent = 36

# --- Paragraph 5 (was "The third and fourth errors ..."): renumber to "The second and third errors ..." ---
Invoke-FindReplace $d.Content `
    "The third and fourth errors (CID 10281 and 10280) were both of the type" `
    "The second and third errors (CID 10281 and 10280) were both of the type" | Out-Null

Write-Output "Edit complete."
